$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

$ws.Range("A2").Value = 19330051920163
$ws.Range("B2").Value = "LOPEZ"
$ws.Range("C2").Value = "SANCHEZ"
$ws.Range("D2").Value = "CINTHIA"
$ws.Range("E2").Value = "INGLÉS V"
$ws.Range("F2").Value = "5ALCM"
$ws.Range("G2").Value = 6

$ws.Range("A3").Value = 19330051920140
$ws.Range("B3").Value = "SOLIS"
$ws.Range("C3").Value = "MARCELINO"
$ws.Range("D3").Value = "MARLENE ALICIA"
$ws.Range("E3").Value = "INGLÉS V"
$ws.Range("F3").Value = "5ARHM"
$ws.Range("G3").Value = 6

$ws.Range("A4").Value = 19330051920195
$ws.Range("B4").Value = "GUILLEN"
$ws.Range("C4").Value = "LINARES"
$ws.Range("D4").Value = "EDITH"
$ws.Range("E4").Value = "INGLÉS V"
$ws.Range("F4").Value = "5BLCM"
$ws.Range("G4").Value = 6

$ws.Range("A5").Value = 19330051920202
$ws.Range("B5").Value = "LEON"
$ws.Range("C5").Value = "GONZALEZ"
$ws.Range("D5").Value = "DANNA PAOLA"
$ws.Range("E5").Value = "INGLÉS V"
$ws.Range("F5").Value = "5BLCM"
$ws.Range("G5").Value = 6

$ws.Range("A6").Value = 18330061460390
$ws.Range("B6").Value = "ROMAN"
$ws.Range("C6").Value = "CASTILLO"
$ws.Range("D6").Value = "NATANAEL"
$ws.Range("E6").Value = "INGLÉS V"
$ws.Range("F6").Value = "5AEM"
$ws.Range("G6").Value = 7

$ws.Range("A7").Value = 19330051920038
$ws.Range("B7").Value = "VAZQUEZ"
$ws.Range("C7").Value = "CHICO"
$ws.Range("D7").Value = "ERICK ORLANDO"
$ws.Range("E7").Value = "INGLÉS V"
$ws.Range("F7").Value = "5AEM"
$ws.Range("G7").Value = 7

$ws.Range("A8").Value = 19330051920276
$ws.Range("B8").Value = "COUDER"
$ws.Range("C8").Value = "SANCHEZ"
$ws.Range("D8").Value = "YULIANA"
$ws.Range("E8").Value = "INGLÉS V"
$ws.Range("F8").Value = "5ALCM"
$ws.Range("G8").Value = 6

$ws.Range("A9").Value = 19330051920102
$ws.Range("B9").Value = "GONZALEZ"
$ws.Range("C9").Value = "IXMATLAHUA"
$ws.Range("D9").Value = "MIGUEL ANGEL"
$ws.Range("E9").Value = "INGLÉS V"
$ws.Range("F9").Value = "5BEM"
$ws.Range("G9").Value = 7

$ws.Range("A10").Value = 19330051920197
$ws.Range("B10").Value = "HERNANDEZ"
$ws.Range("C10").Value = "GUTIERREZ"
$ws.Range("D10").Value = "KEVIN JETHZAEL"
$ws.Range("E10").Value = "INGLÉS V"
$ws.Range("F10").Value = "5BLCM"
$ws.Range("G10").Value = 6

$ws.Range("A11").Value = 19330051920213
$ws.Range("B11").Value = "PEREZ"
$ws.Range("C11").Value = "RAMIREZ"
$ws.Range("D11").Value = "JANETH"
$ws.Range("E11").Value = "INGLÉS V"
$ws.Range("F11").Value = "5BLCM"
$ws.Range("G11").Value = 6

